$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# "Change from parameters to registers and endpoints":
# fold the DAC80508_CONFIG sheet's register rows into the DAC80508 sheet
# (right above the existing WB_IN/WB_OUT endpoint rows), then remove the
# now-redundant DAC80508_CONFIG sheet.
$src = $wb.Worksheets.Item("DAC80508_CONFIG")
$dst = $wb.Worksheets.Item("DAC80508")

# Make room for the 14 register rows (src rows 2:15) right before the
# current row 18 (WB_IN_0 ...).
$dst.Rows("18:31").Insert()

$srcRange = $src.Range("A2:F15")
$dstRange = $dst.Range("A18:F31")
$dstRange.Value2 = $srcRange.Value2

# The register rows now live on DAC80508, so drop the old config sheet.
$src.Delete()

# Leave DAC80508 scrolled down to the newly appended rows.
$dst.Activate()
$dst.Range("A38").Select()

# Active tab moves on to AD5453 (now immediately after DAC80508).
$ad5453 = $wb.Worksheets.Item("AD5453")
$ad5453.Activate()
$ad5453.Range("A2").Select()
